$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update C2: Maximum Voltage (mV) Diseased State value
$ws.Range("C2").Value = 185.21983161833489

# Update C5, C6, C7 to "N/A" (Diseased State not available for these intervals)
$ws.Range("C5").Value = "N/A"
$ws.Range("C6").Value = "N/A"
$ws.Range("C7").Value = "N/A"
